# Add a new "Player Info" sheet as the first sheet in the workbook, and
# replace the MATCH_CARD_LINK (full scorecard URL) column with a
# MATCH_CODE column (just the numeric match code) on both existing sheets.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Write a value as plain text, even when it looks like a number
    # (matches this workbook's convention of storing everything as text),
    # then drop back to the default "Normal" style so no stray number
    # format / style index is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- 1. Insert new "Player Info" sheet before the first existing sheet ---
$playerInfo = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$playerInfo.Name = "Player Info"

# Match the page margins used throughout the rest of the workbook
$playerInfo.PageSetup.LeftMargin = 0.75 * 72
$playerInfo.PageSetup.RightMargin = 0.75 * 72
$playerInfo.PageSetup.TopMargin = 1 * 72
$playerInfo.PageSetup.BottomMargin = 1 * 72
$playerInfo.PageSetup.HeaderMargin = 0.5 * 72
$playerInfo.PageSetup.FooterMargin = 0.5 * 72

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the header style used on the other sheets (bold, centered, bordered)
$srcHeader = $wb.Worksheets.Item("ODI Batting").Range("A1")
$srcHeader.Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Set-TextValue $playerInfo.Range("A2") "7117"
$playerInfo.Range("B2").Value = "Usama Mir"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

# --- 2. Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"
Set-TextValue $battingSheet.Range("D2") "4686"
Set-TextValue $battingSheet.Range("D3") "4688"
Set-TextValue $battingSheet.Range("D4") "4690"

# --- 3. Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
Set-TextValue $bowlingSheet.Range("B2") "4686"
Set-TextValue $bowlingSheet.Range("B3") "4688"
Set-TextValue $bowlingSheet.Range("B4") "4690"
